# [Kadastro App] Kayıtlar güncellendi - 16.07.2025 23:02:42
#
# Adds a header row + one data record to the "Kayitlar" sheet, and turns
# off right-to-left sheet display (was explicitly "0"/false already, the
# authoring tool simply stopped emitting the attribute) on every sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Kayitlar")

$headers = @("Kayıt No", "Tarih", "Birim", "Dosya Sayısı", "Parsel Sayısı", "İş", "Personeller")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$row2 = @("2", "2025-07-16", "Merkez", "3", "2", "Cins D.", "Göktan ELGÜL")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $cell = $ws.Cells.Item(2, $i + 1)
    # Force text storage (incl. the numeric-looking / date-looking values)
    # by entering them with a leading quote, then strip the resulting
    # "quote prefix" style back off so no extra formatting is left behind.
    $cell.Value = "'" + $row2[$i]
    $cell.Style = "Normal"
}

foreach ($sheet in $wb.Worksheets) {
    $sheet.DisplayRightToLeft = $false
}
